# WineCompany.xlsx edit:
#  - Remove the (empty) "Resultados_ElaboracionVino" sheet
#  - Update the "Elaboracion_Vino" quantities
#  - Update the "Distribucion_Ordenes" order-type column
#  - Re-point the active sheet/selection to match the author's new focus

$wb = $excel.ActiveWorkbook

# 1. Delete the empty "Resultados_ElaboracionVino" sheet entirely.
$wsResultados = $wb.Worksheets.Item("Resultados_ElaboracionVino")
$wsResultados.Delete() | Out-Null

# 2. Update "Elaboracion_Vino" quantities (column B).
$wsElab = $wb.Worksheets.Item("Elaboracion_Vino")
$wsElab.Range("B2").Value = 320
$wsElab.Range("B3").Value = 3200
$wsElab.Range("B4").Value = 480
$wsElab.Activate() | Out-Null
$wsElab.Range("B2").Select() | Out-Null

# 3. Update "Distribucion_Ordenes" order-type column (column B, rows 2-14).
$wsDist = $wb.Worksheets.Item("Distribucion_Ordenes")
$newOrderTypes = @(1, 1, 1, 1, 1, 1, 2, 2, 2, 2, 2, 2, 2)
for ($i = 0; $i -lt $newOrderTypes.Length; $i++) {
    $row = 2 + $i
    $wsDist.Cells.Item($row, 2).Value = $newOrderTypes[$i]
}

# 4. "Ventas" sheet keeps its data but its selection moves to F11 and it is
#    no longer the active tab.
$wsVentas = $wb.Worksheets.Item("Ventas")
$wsVentas.Activate() | Out-Null
$wsVentas.Range("F11").Select() | Out-Null

# 5. "Distribucion_Ordenes" becomes the active sheet/tab, selection at B15.
$wsDist.Activate() | Out-Null
$wsDist.Range("B15").Select() | Out-Null
